# Generate Report for Handoff
# Adds two new source files (6ea86943-... and 8cd68c7b-...) to the
# localization status report: one new row per file on the "Overview"
# sheet, and one new row per file on each of the "zh-cn" and "de-de"
# per-locale sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# New file identifiers driving this handoff batch
# ---------------------------------------------------------------
$file1Guid = "6ea86943-b7e0-4f8f-8f2a-891b6309fb6a"
$file2Guid = "8cd68c7b-8843-4ec8-bb6a-1519158efe00"

$file1Name = "$file1Guid.md"
$file2Name = "$file2Guid.md"

$file1Path = "e2e\$file1Name"
$file2Path = "e2e\$file2Name"

$mdDate = "2016-08-25 12:42:45"

$file1ZhXlf = "$file1Guid.f35575f34d2a966988c3460b78ce4689d5c30113.zh-cn.xlf"
$file2ZhXlf = "$file2Guid.e78b1cc9c9ff3b089e1c7c5ab350849947932a89.zh-cn.xlf"
$zhHandoffDate = "2016-08-25 12:42:41"

$file1DeXlf = "$file1Guid.f35575f34d2a966988c3460b78ce4689d5c30113.de-de.xlf"
$file2DeXlf = "$file2Guid.e78b1cc9c9ff3b089e1c7c5ab350849947932a89.de-de.xlf"
$deHandoffDate = "2016-08-25 12:42:45"

$epoch = "0001-01-01 00:00:00"

# =================================================================
# Overview sheet
# =================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A4").Value = $file1Name
$wsOverview.Range("B4").Value = $file1Path
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = $mdDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file1Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file1Path) | Out-Null
$wsOverview.Range("B4").Style = "HyperLink"

$loOverview.ListRows.Add() | Out-Null
$wsOverview.Range("A5").Value = $file2Name
$wsOverview.Range("B5").Value = $file2Path
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = $mdDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file2Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file2Path) | Out-Null
$wsOverview.Range("B5").Style = "HyperLink"

# =================================================================
# zh-cn sheet
# =================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A4").Value = $file1Name
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = $file1ZhXlf
$wsZh.Range("H4").Value = $zhHandoffDate
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $epoch
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file1Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file1Name) | Out-Null
$wsZh.Range("A4").Style = "HyperLink"

$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A5").Value = $file2Name
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = $file2ZhXlf
$wsZh.Range("H5").Value = $zhHandoffDate
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $epoch
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "False"
$wsZh.Range("P5").Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file2Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file2Name) | Out-Null
$wsZh.Range("A5").Style = "HyperLink"

# =================================================================
# de-de sheet
# =================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A4").Value = $file1Name
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = $file1DeXlf
$wsDe.Range("H4").Value = $deHandoffDate
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $epoch
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file1Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file1Name) | Out-Null
$wsDe.Range("A4").Style = "HyperLink"

$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A5").Value = $file2Name
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = $file2DeXlf
$wsDe.Range("H5").Value = $deHandoffDate
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $epoch
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "False"
$wsDe.Range("P5").Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/$file2Path".Replace('\','/'), [System.Type]::Missing, [System.Type]::Missing, $file2Name) | Out-Null
$wsDe.Range("A5").Style = "HyperLink"

Write-Host "Handback report updated with $file1Name and $file2Name"
